# projectSapper.pptx — "Add files via upload" content fixes
#
# Slide 2 ("ЦЕЛИ ПРОЕКТА"): the two goal bullets get reflowed into several
# same-formatted runs (no visible text change, just a run split that
# happened when the author re-typed the paragraphs).
#
# Slide 5 ("Проблемы"): the CISC/RISC bullets are likewise reflowed into
# several same-formatted runs, and two of the bullets also get their
# wording corrected:
#   "RISC использует статичную длинну операции"   -> "RISC использует статичную длину операций"
#   "CISC использует статичную длинну операции"   -> "CISC использует динамическую длину операций"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 — PlaceHolder 2 (goals bullet list)
# ---------------------------------------------------------------------
$s2  = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange

# Paragraph 1: "Научиться разрабатывать кроссплатформенные и кроссархитектурные приложения  "
$para = $tr2.Paragraphs(1, 1)
$para.Characters(1, 24).Text  = "Научиться разрабатывать "
$para.Characters(25, 21).Text = "кроссплатформенные и "
$para.Characters(46, 19).Text = "кроссархитектурные "
$para.Characters(65, 12).Text = "приложения  "

# Paragraph 2: "Разработать игру "Сапер" под Linux и Windows"
$para = $tr2.Paragraphs(2, 1)
$para.Characters(1, 25).Text  = "Разработать игру “Сапер” "
$para.Characters(26, 19).Text = "под Linux и Windows"

# ---------------------------------------------------------------------
# Slide 5 — PlaceHolder 3 (CISC vs RISC bullet list)
# ---------------------------------------------------------------------
$s5  = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(3)
$tr5 = $sh5.TextFrame.TextRange

# Paragraph 2: "CISC поддержка многообразных сложных инструкций выполняющихся более чем за один такт"
$para = $tr5.Paragraphs(2, 1)
$para.Characters(1, 37).Text  = "CISC поддержка многообразных сложных "
$para.Characters(38, 43).Text = "инструкций выполняющихся более чем за один "
$para.Characters(81, 4).Text  = "такт"

# Paragraph 3: "RISC минимум инструкций, стремящихся к выполнению за один такт"
$para = $tr5.Paragraphs(3, 1)
$para.Characters(1, 39).Text  = "RISC минимум инструкций, стремящихся к "
$para.Characters(40, 23).Text = "выполнению за один такт"

# Paragraph 4: wording fix only (no run split)
#   "RISC использует статичную длинну операции" -> "RISC использует статичную длину операций"
$para = $tr5.Paragraphs(4, 1)
$len = $para.Text.Length
$para.Characters(1, $len).Text = "RISC использует статичную длину операций"

# Paragraph 5: wording fix + run split
#   "CISC использует статичную длинну операции" -> "CISC использует динамическую длину операций"
$para = $tr5.Paragraphs(5, 1)
$len = $para.Text.Length
$para.Characters(1, $len).Text = "CISC использует динамическую длину операций"
$para.Characters(1, 35).Text  = "CISC использует динамическую длину "
$para.Characters(36, 8).Text  = "операций"
